$d = $word.ActiveDocument

# 1. Title: "FACT December Workshops 2020" -> "GLATOS April Workshops 2021"
$d.Content.Find.Execute("FACT December Workshops 2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GLATOS April Workshops 2021", 2)

# 2. Insert new paragraph "Detailed installation instructions are below." after the
#    WARNING paragraph (it inherits the WARNING paragraph's indentation formatting).
$warningPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "WARNING: Software installation*") {
        $warningPara = $p
        break
    }
}
$warningPara.Range.InsertParagraphAfter()
$newParaIndex = $warningPara.Index + 1
$newPara = $d.Paragraphs($newParaIndex)
$newPara.Range.Text = "Detailed installation instructions are below."

# 3. "PLEASE ONLY DOWNLOAD 2020-12-15 (we are still working on the content):" ->
#    "PLEASE ONLY DOWNLOAD THE MORNING OF WORKSHOP (we are still working on the content)"
#    (drop the trailing colon, change the date text to "THE MORNING OF WORKSHOP")
$d.Content.Find.Execute("EASE ONLY DOWNLOAD 2020-12-15 (we are still working on the content):", $true, $false, $false, $false, $false,
                         $true, 1, $false, "EASE ONLY DOWNLOAD THE MORNING OF WORKSHOP (we are still working on the content)", 2)

# 4. GitHub dataset/workshop repo link + display text
foreach ($h in $d.Hyperlinks) {
    if ($h.Address -eq "https://github.com/ocean-tracking-network/2020-12-17-telemetry-packages-FACT/") {
        $h.TextToDisplay = "https://github.com/ocean-tracking-network/2021-03-30-glatos-workshop"
        $h.Address = "https://github.com/ocean-tracking-network/2021-03-30-glatos-workshop"
        break
    }
}

# 5. Add trailing period to the "Select the GREEN "code" button..." paragraph
$d.Content.Find.Execute('Select the GREEN "code" button at the top and choose "Download ZIP"', $true, $false, $false, $false, $false,
                         $true, 1, $false, 'Select the GREEN "code" button at the top and choose "Download ZIP".', 2)

# 6. git clone URL update
$d.Content.Find.Execute("https://github.com/ocean-tracking-network/2020-12-17-telemetry-packages-FACT.git", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://github.com/ocean-tracking-network/2021-03-30-glatos-workshop.git", 2)

# 7. "This is only needed for Windows users - not MacOS" -> "Only Windows Users need to install Rtools - not MacOS"
$d.Content.Find.Execute("This is only needed for Windows users" + [char]8211 + " not MacOS", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Only Windows Users need to install Rtools" + [char]8211 + " not MacOS", 2)

"done"
